$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Bestand 2" values in column C with "Bestand 1" values, and
# highlight the changed cells with a solid red fill (matching the other
# workbook's conflicting data so it stands out).
$updates = @{
    "C4"  = "Bestand 1 column 4"
    "C6"  = "Bestand 1 column 6"
    "C7"  = "Bestand 1 column 2"
    "C9"  = "Bestand 1 column 8"
    "C10" = "Bestand 1 column 11"
    "C11" = "Bestand 1 column 10"
    "C12" = "Bestand 1 column 9"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Interior.Color = 255
}
